$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 - this shifts existing rows 2..16 down to 3..17
# and Excel auto-adjusts the formula references (B6->B7, B9->B10, shared-formula
# ranges, etc.) exactly like a native Insert operation.
$ws.Rows("2:2").Insert()

# The newly inserted row comes back blank/unformatted - carry the date
# number format down from the row beneath it (mirrors what Excel does when a
# user fills the new row in by hand after inserting above the top data row).
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# New row 2: latest data point
$ws.Range("A2").Value = 45865
$ws.Range("B2").Value = 46743
$ws.Range("C2").Formula = "=(B2/B7-1)*100"

# The two rows that used to be row 2 and row 3 (now row 3 and row 4) get
# their date refreshed to the newer reporting date while keeping their
# existing Value (B) figures / formulas.
$ws.Range("A3").Value = 45774
$ws.Range("A4").Value = 45683

# Row 3's comparison formula was manually re-pointed one row further down
# (B7 -> B8) instead of being left at the value the insert auto-shifted it to.
$ws.Range("C3").Formula = "=(B3/B8-1)*100"

# Selection left on B3 by the editor at save time.
$ws.Range("B3").Select()
